# Add new "welcome screens" translation rows to the translations sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# key, en, hi, te
$rows = @(
    @("skip", "Skip", "छोड़ें", "దాటవేయి"),
    @("title", "Find Perfect Properties", "परफेक्ट प्रॉपर्टी खोजें", "పరిపూర్ణ ఆస్తులను కనుగొనండి"),
    @("next", "Next", "आगे", "తదుపరి"),
    @("smart_search_label", "Smart Search", "स्मार्ट खोज", "స్మార్ట్ శోధన"),
    @("direct_property_connect_label", "Direct Property Connect", "प्रत्यक्ष प्रॉपर्टी कनेक्ट", "ప్రత్యక్ష ప్రాపర్టీ కనెక్ట్"),
    @("discover_title", "Discover a place you will love to live", "एक ऐसी जगह खोजें जहाँ आप रहना पसंद करेंगे", "మీరు నివసించడం ఇష్టపడే స్థలాన్ని కనుగొనండి"),
    @("discover_description", "Browse homes for sale, original neighborhood photos, resident reviews, local insights to find what is right for you", "बिक्री के लिए घर देखें, मूल पड़ोस की तस्वीरें, निवासी समीक्षाएं और स्थानीय जानकारी से अपने लिए सही जगह पाएं", "అమ్మకానికి ఇళ్లను, ప్రాంతీయ ఫోటోలు, నివాసుల సమీక్షలు మరియు లోకల్ సమాచారంతో మీకు సరైన ప్రదేశాన్ని కనుగొనండి"),
    @("lets_explore_button", "Let's Explore", "चलिए शुरू करें", "చదుకుదాం / ముందుకు సాగుదాం")
)

$startRow = 393
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
    $ws.Cells.Item($r, 4).Value = $rows[$i][3]
}

# Update view state to match where the author ended up scrolled to/selected.
$excel.ActiveWindow.ScrollRow = 380
$ws.Range("A403").Select()
